{"js": "// Apply resume text edits described by the diff:\n//  1) Soft Skills paragraph: capitalize each skill keyword, and fix\n//     \"ability to read\" -> \"Ability to read\" (Fast/Friendly too).\n//  2) \"Dono Application\" Role paragraph: rewrite the MongoDB/Mongoose\n//     sentence.\n//  3) \"Workout Buddies\" Role paragraph: rewrite the MySQL/Sequelize /\n//     client-server sentence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Helper: replace the first (and only expected) case-sensitive match of\n// `searchText` inside a given paragraph with `replaceText`.\nasync function replaceInParagraph(paragraph, searchText, replaceText) {\n  const results = paragraph.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Locate the three target paragraphs by distinctive text fragments so the\n// script is resilient to exact paragraph index changes.\nlet softSkillsPara = null;\nlet donoRolePara = null;\nlet workoutRolePara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"leadership, communication, honest\") !== -1) {\n    softSkillsPara = paragraphs.items[i];\n  } else if (t.indexOf(\"on controllers manage MongoDB\") !== -1) {\n    donoRolePara = paragraphs.items[i];\n  } else if (t.indexOf(\"Created and manage MySQL\") !== -1) {\n    workoutRolePara = paragraphs.items[i];\n  }\n}\n\nif (!softSkillsPara) throw new Error(\"Soft Skills paragraph not found\");\nif (!donoRolePara) throw new Error(\"Dono Application role paragraph not found\");\nif (!workoutRolePara) throw new Error(\"Workout Buddies role paragraph not found\");\n\n// 1) Soft Skills line -------------------------------------------------\nawait replaceInParagraph(\n  softSkillsPara,\n  \"leadership, communication, honest, humble, hardworking, critical Thinker, problem solver\",\n  \"Leadership, Communication, Honest, Humble, Hardworking, Critical Thinker, Problem Solver\"\n);\nawait replaceInParagraph(\n  softSkillsPara,\n  \"fast learner, friendly, ability to read and understand documentation\",\n  \"Fast learner, Friendly, Ability to read and understand documentation\"\n);\n\n// 2) Dono Application \"Role:\" line ------------------------------------\nawait replaceInParagraph(\n  donoRolePara,\n  \"Created controllers and routes, on controllers manage MongoDB (Mongoose) collections using CRUD and filter response based on request. Used React.Js global context to store server response to later render data to components. Also created Passport JWT authentication, refactor components, created custom hooks, \",\n  \"Created controllers and routes; manage MongoDB collections using Mongoose and filter or populated response based on request. Used React.Js global context to store server response to later render data. Also created Passport JWT authentication, refactor components, created custom hooks. \"\n);\n\n// 3) Workout Buddies \"Role:\" line --------------------------------------\nawait replaceInParagraph(\n  workoutRolePara,\n  \"Created and manage MySQL (Sequelize), used routes to manage database and filter data based on client request, get routes used server side rendering with handlebars.js, on client side used bootstrap, JQuery, client-side API calls to server with AJAX, and server-side RESTful API request with axios. Also created Passport Sessions authentication.\",\n  \"Created and manage MySQL database with Sequelize, used routes to manage database and filter data based on client request, get routes used server side rendering with handlebars.js. Client side used bootstrap, JQuery, API calls to server with AJAX. Server-side used RESTful API request with axios. Also created Passport Sessions authentication.\"\n);\n", "ps1": "# Apply resume text edits described by the diff:\n#  1) Soft Skills paragraph: capitalize each skill keyword, and fix\n#     \"ability to read\" -> \"Ability to read\" (Fast/Friendly too).\n#  2) \"Dono Application\" Role paragraph: rewrite the MongoDB/Mongoose\n#     sentence.\n#  3) \"Workout Buddies\" Role paragraph: rewrite the MySQL/Sequelize /\n#     client-server sentence.\n\n$d = $word.ActiveDocument\n\n# Locate the three target paragraphs by distinctive text fragments so the\n# script is resilient to exact paragraph index changes.\n$softSkillsIndex = -1\n$donoRoleIndex = -1\n$workoutRoleIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*leadership, communication, honest*\") {\n        $softSkillsIndex = $i\n    } elseif ($t -like \"*on controllers manage MongoDB*\") {\n        $donoRoleIndex = $i\n    } elseif ($t -like \"*Created and manage MySQL*\") {\n        $workoutRoleIndex = $i\n    }\n}\n\nif ($softSkillsIndex -eq -1) { throw \"Soft Skills paragraph not found\" }\nif ($donoRoleIndex -eq -1) { throw \"Dono Application role paragraph not found\" }\nif ($workoutRoleIndex -eq -1) { throw \"Workout Buddies role paragraph not found\" }\n\nfunction Replace-InParagraph([int]$paraIndex, [string]$searchText, [string]$replaceText) {\n    $p = $d.Paragraphs.Item($paraIndex)\n    $r = $p.Range\n    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $searchText\"\n    }\n}\n\n# 1) Soft Skills line -------------------------------------------------\nReplace-InParagraph $softSkillsIndex `\n    \"leadership, communication, honest, humble, hardworking, critical Thinker, problem solver\" `\n    \"Leadership, Communication, Honest, Humble, Hardworking, Critical Thinker, Problem Solver\"\n\nReplace-InParagraph $softSkillsIndex `\n    \"fast learner, friendly, ability to read and understand documentation\" `\n    \"Fast learner, Friendly, Ability to read and understand documentation\"\n\n# 2) Dono Application \"Role:\" line ------------------------------------\nReplace-InParagraph $donoRoleIndex `\n    \"Created controllers and routes, on controllers manage MongoDB (Mongoose) collections using CRUD and filter response based on request. Used React.Js global context to store server response to later render data to components. Also created Passport JWT authentication, refactor components, created custom hooks, \" `\n    \"Created controllers and routes; manage MongoDB collections using Mongoose and filter or populated response based on request. Used React.Js global context to store server response to later render data. Also created Passport JWT authentication, refactor components, created custom hooks. \"\n\n# 3) Workout Buddies \"Role:\" line --------------------------------------\nReplace-InParagraph $workoutRoleIndex `\n    \"Created and manage MySQL (Sequelize), used routes to manage database and filter data based on client request, get routes used server side rendering with handlebars.js, on client side used bootstrap, JQuery, client-side API calls to server with AJAX, and server-side RESTful API request with axios. Also created Passport Sessions authentication.\" `\n    \"Created and manage MySQL database with Sequelize, used routes to manage database and filter data based on client request, get routes used server side rendering with handlebars.js. Client side used bootstrap, JQuery, API calls to server with AJAX. Server-side used RESTful API request with axios. Also created Passport Sessions authentication.\"\n"}
